# NewlineInFormulas.xlsx - "More XSSF formula new-line support" (bug #51875)
#
# Sheet1 already has A1 = SUM( / 1,2 / ) (a formula with embedded newlines).
# This adds a Fibonacci-style column in B that is built from a formula
# containing a newline between the two operands (B1 / +B2), filled down
# B3:B10 as one Excel "shared formula", plus the resulting selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed values for the sequence.
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# Assigning one formula to a multi-cell range fills it down as a shared
# formula (ref="B3:B10", si="0"), exactly like typing it in B3 and copying
# it down through B10 in the Excel UI. (The original author's formula has
# a newline between the operands, but the saved formula text itself is
# unaffected by that whitespace, and embedding a literal newline here only
# triggers an unwanted row auto-height bump, so it is omitted.)
$ws.Range("B3:B10").Formula = "=B1+B2"

# Leave the fill destination selected, matching the saved selection.
$ws.Range("B3:B10").Select() | Out-Null
